$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns with the latest
# scraped values from the GitHub Actions cron job.
# A leading apostrophe is used for Price values that would otherwise
# be auto-recognized by Excel as numbers, forcing them to remain text
# (matching the original inline-string cell contents).

$ws.Range("D2").Value = "68.931.00"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "2.465.60"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'557.87"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "'162.61"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").Value = "2.464.55"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -6.07%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "'0.334"
$ws.Range("E12").Value = "  -4.91%  "
$ws.Range("D13").Value = "'4.82"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "2.917.31"
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "68.796.08"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("D17").Value = "'23.54"
$ws.Range("E17").Value = "  -2.95%  "
$ws.Range("D18").Value = "2.466.46"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  -4.50%  "
$ws.Range("D20").Value = "'341.92"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").Value = "'7.04"
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").Value = "'3.79"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").Value = "'66.87"
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").Value = "2.592.01"
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "'8.13"
$ws.Range("E30").Value = "  -5.73%  "
$ws.Range("D31").Value = "0.0₃0819"
$ws.Range("E31").Value = "  -6.40%  "
$ws.Range("D32").Value = "'7.18"
$ws.Range("E32").Value = "  -5.96%  "
$ws.Range("D33").Value = "'437.98"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -4.56%  "
$ws.Range("E36").Value = "  -5.87%  "
$ws.Range("D37").Value = "'157.37"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("D38").Value = "'19.06"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("D41").Value = "'17.86"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("D42").Value = "'0.303"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("D44").Value = "'37.44"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("E45").Value = "  -6.67%  "
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("E47").Value = "  -5.50%  "
$ws.Range("D48").Value = "'132.82"
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").Value = "'0.0715"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "'0.484"
$ws.Range("E51").Value = "  -4.42%  "
